$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reset the task-table region to a clean slate ---
$ws.Range("A1:H22").Style = "Normal"
$ws.Range("A1:H22").ClearContents()

# --- Write cell values ---
$ws.Cells.Item(1,1).Value = 'S.No.'
$ws.Cells.Item(1,2).Value = 'Task'
$ws.Cells.Item(1,3).Value = 'Owner'
$ws.Cells.Item(1,4).Value = 'Start Date'
$ws.Cells.Item(1,5).Value = 'End Date'
$ws.Cells.Item(1,6).Value = 'Status'
$ws.Cells.Item(1,7).Value = 'Comments'
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 'Create addon "storeActivityAddOn" and install it on bncwebservices'
$ws.Cells.Item(2,3).Value = 'Swarnima/Swapnil'
$ws.Cells.Item(2,4).Value = '23/03'
$ws.Cells.Item(2,5).Value = '23/03'
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 'Create an item "storeActivity" and add following attributes to it.'+[char]10+'1. storeId'+[char]10+'2. customerId'+[char]10+'3. storeVisitDate'+[char]10+'4. storeEntryTime'+[char]10+'5.storeExitTime'+[char]10+'6. timeSpentInStore (save time in minutes, to be calculated only when storeEntryTime and storeExitTime are available)'
$ws.Cells.Item(3,3).Value = 'Swapnil'
$ws.Cells.Item(3,4).Value = '24/03'
$ws.Cells.Item(3,5).Value = '24/03'
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 'Create a web service which will return the status (status:entered/exited) in JSON format'+[char]10+'URL(/storeActivity)'+[char]10+'It will receive the following parameters:'+[char]10+'1. storeId'+[char]10+'2. customerId'+[char]10+'3. storeVisitDate'
$ws.Cells.Item(4,7).Value = 'First there will be a check whether the customer''s entry exists, in case it exists and exit time is null in that, then it will calculate the timeSpentInStore and update the same along with the exit time.(as new Date())'+[char]10+'In case there is an entry and exit time is not null, also if there is no entry for the customer, then it will create a new entry and save the same along with the entry time (as new Date())'+[char]10+'calculate the following field in the controller method'+[char]10+'timeSpentInStore (save time in minutes, to be calculated only when storeEntryTime and storeExitTime are available)'
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 'Create dao, service classes for the above point'
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 'Create some dummy data for storeActivity itemtype(Impex generation also)'
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 'Create service and dao classes to get most visited stores by customers. The "storeActivity" item (created in step 2) is to be queried to get most visited stores'
$ws.Cells.Item(7,3).Value = 'Swarnima'
$ws.Cells.Item(7,4).Value = '24/03'
$ws.Cells.Item(7,5).Value = '25/03'
$ws.Cells.Item(7,7).Value = 'Property file - for showing a 5-6 most visited stores'
$ws.Cells.Item(7,8).Value = 'Impex'
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 'Create service and dao classes to get loyal customers. The "storeActivity" item is to be queried to get the customers who visit stores the most'
$ws.Cells.Item(8,3).Value = 'Swarnima'
$ws.Cells.Item(8,4).Value = '26/03'
$ws.Cells.Item(8,5).Value = '27/03'
$ws.Cells.Item(8,7).Value = 'Property file - for showing a 5-6 most loyal customers'
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 'Create an item type "beacon" that contains following String type attributes.'+[char]10+'1. beaconId'+[char]10+'2. majorId'+[char]10+'3. minorId'
$ws.Cells.Item(9,3).Value = 'Swapnil'
$ws.Cells.Item(9,4).Value = '26/03'
$ws.Cells.Item(9,5).Value = '26/03'
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 'Create an enumtype BeaconType with following Types'+[char]10+'1. Entry'+[char]10+'2. Exit'+[char]10+'3. Product'+[char]10+'4. Checkout'
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 'Create a relation between beacon and BeaconType, a one to many relation(create impex add some dummy data for beacons and relation)'
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = 'Customize product item type by adding following attribute.'+[char]10+'1. popularityCount (int type)'
$ws.Cells.Item(12,3).Value = 'Swapnil'
$ws.Cells.Item(12,4).Value = '26/03'
$ws.Cells.Item(12,5).Value = '26/03'
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = 'Create a one to many relation between beacon and products. And add some sample data in it (create impex)'
$ws.Cells.Item(13,3).Value = 'Swapnil'
$ws.Cells.Item(13,4).Value = '26/03'
$ws.Cells.Item(13,5).Value = '26/03'
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 'Create a webservice to update the popularityCount of the product. The web service will accept the beaconId, majorId, minorId (in JSON format).'+[char]10+'URL(/popularityCount)'+[char]10+' Will use the request data to get the product from the relation defined in point 8 above and increase the popularity count for the product.'
$ws.Cells.Item(14,3).Value = 'Swapnil'
$ws.Cells.Item(14,4).Value = '26/03'
$ws.Cells.Item(14,5).Value = '31/03'
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = 'Create service and dao classes to get top 5 popular products based on the popularity count. The count "5" should be made configurable.'
$ws.Cells.Item(15,3).Value = 'Swapnil'
$ws.Cells.Item(15,4).Value = 42095
$ws.Cells.Item(15,5).Value = 42096
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = 'Create a webservice to return the beacon type in JSON format. The web service will accept the beaconId, majorId, minorId (in JSON format).'+[char]10+'URL(/beaconType)'+[char]10+'And it will query and return the type of the beacon'
$ws.Cells.Item(16,7).Value = 'At the IOS side, we will need to have a link such that, when entry/exit beacon type is found, then they should send the webservice with URL /storeActivity'+[char]10+'When product beacon type is found then they should send the webservice with URL /popularityCount'
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = 'Create service and dao classes for the above point'
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = 'create an item "storeCustomer" and add following attributes to it.'+[char]10+'1. weight'+[char]10+'2. height'+[char]10+'3. Age'+[char]10+'4. Gender'+[char]10+''
$ws.Cells.Item(18,3).Value = 'Swapnil'
$ws.Cells.Item(18,4).Value = 42097
$ws.Cells.Item(18,5).Value = 42097
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = 'Create impex to store some sample data in storeCustomer point 11 above.'
$ws.Cells.Item(19,3).Value = 'Swapnil'
$ws.Cells.Item(19,4).Value = 42097
$ws.Cells.Item(19,5).Value = 42097
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = 'Create a WCMS page for the Activity dashboard, impex creation.'
$ws.Cells.Item(20,3).Value = 'Swarnima'
$ws.Cells.Item(20,4).Value = '30/03'
$ws.Cells.Item(20,5).Value = '30/03'
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = 'Create controller, facades, service classes to fetch model data for store customer profile section (point 11), popular products (based on product popularity count, point 7), most visited stores (point 4), loyal customers (from point 5), spent time (query storeActivity item, point 2)'
$ws.Cells.Item(21,3).Value = 'Swarnima'
$ws.Cells.Item(21,4).Value = '31/03'
$ws.Cells.Item(21,5).Value = 42096
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = 'To get weather information, there are two approaches. '+[char]10+'1) The weather data is passed by IOS app, if so we can use it to dispplay on our page.'+[char]10+'2) To use java APIs to fetch weather data based on city or zip code. Please see the sample code http://code.aksingh.net/owm-japis/src'
$ws.Cells.Item(22,3).Value = 'Swarnima'
$ws.Cells.Item(22,4).Value = 42097
$ws.Cells.Item(22,5).Value = 42102

# --- WrapText style (matches original cellXf s=1) ---
$ws.Range("B3,B4,G4,B5,B6,B7,G7,B8,G8,B9,B10,B11,B12,B13,B14,B15,B16,G16,B17,G17,B18,B19,B20,B21,B22").WrapText = $true

# --- Date number format (matches original cellXf s=2, numFmtId 16 'd-mmm') ---
$ws.Range("D15,E15,D16,D17,D18,E18,D19,E19,E21,D22,E22").NumberFormat = "d-mmm"

# --- Header fill (grey, Background 2 Darker 25%) ---
$headerRng = $ws.Range("A1:G1")
$headerRng.Interior.ThemeColor = 4
$headerRng.Interior.TintAndShade = -0.249977111117893

# --- "Good" cell style (green) applied to the now-empty G21 ---
$ws.Range("G21").Style = "Good"

# --- Row heights (explicit, matches target ht values) ---
$ws.Rows.Item(3).RowHeight = 132
$ws.Rows.Item(4).RowHeight = 165.75
$ws.Rows.Item(5).RowHeight = 30.75
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 75
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 90
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 90
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 75
$ws.Rows.Item(22).RowHeight = 75

# --- Column widths ---
$ws.Columns.Item(7).ColumnWidth = 90.6

# --- AutoFilter over C1:C22 (extends the existing filter range) ---
$ws.AutoFilterMode = $false
$ws.Range("C1:C22").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match the new filter range ---
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name() -eq "Sheet1!_FilterDatabase") {
        $nm.RefersTo = "=Sheet1!`$C`$1:`$C`$22"
    }
}

# --- Sheet view: scroll position + selection ---
$ws.Activate()
$excel.Goto($ws.Range("A16"), $true)
$ws.Range("F17").Select()
